# Project "Sample Project" is saved. Update the "Rules" sheet:
# cell B11 (row for the last rule, previously holding the label "R40")
# is changed to contain the text value "1".
#
# NOTE: B11 must keep being a *text* cell (shared string), not become a
# number, and it must keep its existing cell style (border/format) as-is.
# A plain `$ws.Range("B11").Value = "1"` would be auto-coerced by Excel
# into the number 1 (and picking up a different style), so instead we
# write a text-producing formula and then collapse it down to a static
# value via Copy/PasteSpecial(values), which preserves both the text
# type and the original formatting of the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues: paste only the resulting value, keep formatting
$excel.CutCopyMode = $false
